$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 412, shifting existing rows 412:509 down to 413:510.
$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with its data.
$ws.Cells.Item(412, 1).Value = 1
$ws.Cells.Item(412, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(412, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(412, 4).Value = 45211
$ws.Cells.Item(412, 5).Value = 15
$ws.Cells.Item(412, 6).Value = "Fruta"
$ws.Cells.Item(412, 7).Value = 100102
$ws.Cells.Item(412, 8).Value = "Cítricos"
$ws.Cells.Item(412, 9).Value = 100102003
$ws.Cells.Item(412, 10).Value = "Limón"
$ws.Cells.Item(412, 11).Value = "Tahití"
$ws.Cells.Item(412, 12).Value = "Primera"
$ws.Cells.Item(412, 13).Value = 270
$ws.Cells.Item(412, 14).Value = 42000
$ws.Cells.Item(412, 15).Value = 43000
$ws.Cells.Item(412, 16).Value = 42500
$ws.Cells.Item(412, 17).Value = "$/caja 24 kilos"
$ws.Cells.Item(412, 18).Value = "Perú"
$ws.Cells.Item(412, 19).Value = 1771
$ws.Cells.Item(412, 20).Value = 24
